$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "meteorologico"
$ws.Range("B2").Value = "coordinare"
$ws.Range("C2").Value = "trentennale"
$ws.Range("D2").Value = "geografico"
$ws.Range("E2").Value = "coordinata"
$ws.Range("F2").Value = "caldo"
$ws.Range("G2").Value = "climatologico"
$ws.Range("H2").Value = "lacerazione"
$ws.Range("I2").Value = "fabiano"
$ws.Range("J2").Value = "integrazione"

# Row 3
$ws.Range("A3").Value = "sala"
$ws.Range("B3").Value = "pinacoteca"
$ws.Range("C3").Value = "transetto"
$ws.Range("D3").Value = "abside"
$ws.Range("E3").Value = "attribuire"
$ws.Range("F3").Value = "descrizione"
$ws.Range("G3").Value = "culturaitalia"
$ws.Range("H3").Value = "sec"
$ws.Range("I3").Value = "dipinto"
$ws.Range("J3").Value = "francescano"

# Row 4
$ws.Range("A4").Value = "abitante"
$ws.Range("B4").Value = "situare"
$ws.Range("C4").Value = "paese"
$ws.Range("D4").Value = "frazione"
$ws.Range("E4").Value = "castello"
$ws.Range("F4").Value = "località"
$ws.Range("G4").Value = "strada"
$ws.Range("H4").Value = "monte"
$ws.Range("I4").Value = "geografia"
$ws.Range("J4").Value = "territorio"

# Row 5
$ws.Range("B5").Value = "flora"
$ws.Range("C5").Value = "idrografico"
$ws.Range("D5").Value = "gettare"
$ws.Range("F5").Value = "cinghiale"
$ws.Range("G5").Value = "cresta"
$ws.Range("H5").Value = "cascata"
$ws.Range("I5").Value = "profondità"
$ws.Range("J5").Value = "escursionistico"

# Row 6
$ws.Range("A6").Value = "fulvio"
$ws.Range("B6").Value = "censimento"
$ws.Range("C6").Value = "corgna"
$ws.Range("D6").Value = "sagra"
$ws.Range("E6").Value = "toscana"
$ws.Range("F6").Value = "mastioo"
$ws.Range("G6").Value = "ascanio"
$ws.Range("H6").Value = "cortona"
$ws.Range("I6").Value = "anagno"
$ws.Range("J6").Value = "rimuovere"

# Row 7
$ws.Range("A7").Value = "acerbo"
$ws.Range("B7").Value = "lacuno"
$ws.Range("C7").Value = "tessile"
$ws.Range("D7").Value = "interpretabile"
$ws.Range("E7").Value = "calce"
$ws.Range("F7").Value = "alteravare"
$ws.Range("G7").Value = "gigante"
$ws.Range("H7").Value = "ciclopico"
$ws.Range("J7").Value = "indirizzare"

# Row 8
$ws.Range("A8").Value = "con"
$ws.Range("B8").Value = "essere"
$ws.Range("C8").Value = "che"
$ws.Range("D8").Value = "uno"
$ws.Range("E8").Value = "cui"
$ws.Range("F8").Value = "per"
$ws.Range("G8").Value = "più"
$ws.Range("H8").Value = "avere"
$ws.Range("I8").Value = "venire"
$ws.Range("J8").Value = "due"
